$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on Price/Volume columns so numeric-looking strings
# (e.g. "0.9980", "23.380.16") are preserved exactly as text, not coerced to numbers.
$ws.Range("D2:D51").NumberFormat = "@"
$ws.Range("E2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "23.380.16"
$ws.Range("E2").Value = "  +1.59%  "
$ws.Range("D3").Value = "1.627.54"
$ws.Range("E3").Value = "  +2.70%  "
$ws.Range("D4").Value = "0.9980"
$ws.Range("E4").Value = "  -0.58%  "
$ws.Range("D5").Value = "307.12"
$ws.Range("E5").Value = "  +2.06%  "
$ws.Range("D6").Value = "0.9970"
$ws.Range("D7").Value = "0.3786"
$ws.Range("E7").Value = "  +0.67%  "
$ws.Range("D8").Value = "53.08"
$ws.Range("E8").Value = "  +4.77%  "
$ws.Range("D9").Value = "0.3661"
$ws.Range("E9").Value = "  +2.21%  "
$ws.Range("D10").Value = "1.280"
$ws.Range("E10").Value = "  +5.18%  "
$ws.Range("D11").Value = "0.08201"
$ws.Range("E11").Value = "  +1.97%  "
$ws.Range("D12").Value = "0.9984"
$ws.Range("E12").Value = "  -0.53%  "
$ws.Range("D13").Value = "23.21"
$ws.Range("E13").Value = "  +5.82%  "
$ws.Range("D14").Value = "6.670"
$ws.Range("E14").Value = "  +3.28%  "
$ws.Range("D15").Value = "7.467"
$ws.Range("E15").Value = "  +2.11%  "
$ws.Range("D16").Value = "0.00001264"
$ws.Range("E16").Value = "  +3.48%  "
$ws.Range("D17").Value = "1.625.68"
$ws.Range("E17").Value = "  +2.34%  "
$ws.Range("D18").Value = "94.76"
$ws.Range("E18").Value = "  +2.81%  "
$ws.Range("D19").Value = "0.06948"
$ws.Range("E19").Value = "  +2.29%  "
$ws.Range("D20").Value = "18.41"
$ws.Range("E20").Value = "  +2.73%  "
$ws.Range("D21").Value = "6.588"
$ws.Range("E21").Value = "  +2.47%  "
$ws.Range("D22").Value = "0.9979"
$ws.Range("E22").Value = "  -0.51%  "
$ws.Range("E23").Value = "  +1.64%  "
$ws.Range("D24").Value = "23.401.56"
$ws.Range("E24").Value = "  +1.65%  "
$ws.Range("D25").Value = "3.130"
$ws.Range("E25").Value = "  +12.80%  "
$ws.Range("D26").Value = "2.428"
$ws.Range("E26").Value = "  +2.47%  "
$ws.Range("D27").Value = "21.42"
$ws.Range("E27").Value = "  +3.23%  "
$ws.Range("D28").Value = "150.70"
$ws.Range("E28").Value = "  +2.34%  "
$ws.Range("D29").Value = "5.282"
$ws.Range("E29").Value = "  +1.55%  "
$ws.Range("D30").Value = "136.52"
$ws.Range("E30").Value = "  +2.89%  "
$ws.Range("D31").Value = "2.418"
$ws.Range("E31").Value = "  +2.82%  "
$ws.Range("D32").Value = "6.945"
$ws.Range("E32").Value = "  +6.57%  "
$ws.Range("D33").Value = "1.803.72"
$ws.Range("E33").Value = "  +2.09%  "
$ws.Range("D34").Value = "0.9749"
$ws.Range("E34").Value = "  +3.83%  "
$ws.Range("D35").Value = "0.02809"
$ws.Range("E35").Value = "  +5.18%  "
$ws.Range("D36").Value = "10.48"
$ws.Range("E36").Value = "  +4.91%  "
$ws.Range("D37").Value = "0.07473"
$ws.Range("E37").Value = "  +1.67%  "
$ws.Range("D38").Value = "6.248"
$ws.Range("E38").Value = "  +3.15%  "
$ws.Range("D39").Value = "0.2535"
$ws.Range("E39").Value = "  +2.52%  "
$ws.Range("E40").Value = "  +1.00%  "
$ws.Range("E41").Value = "  +5.61%  "
$ws.Range("D42").Value = "0.7172"
$ws.Range("E42").Value = "  +4.43%  "
$ws.Range("E43").Value = "  +7.49%  "
$ws.Range("D44").Value = "16.25"
$ws.Range("D45").Value = "0.6621"
$ws.Range("E45").Value = "  +3.45%  "
$ws.Range("E46").Value = "  +5.33%  "
$ws.Range("D47").Value = "4.032"
$ws.Range("E47").Value = "  +1.09%  "
$ws.Range("D48").Value = "0.9958"
$ws.Range("E48").Value = "  -0.59%  "
$ws.Range("D49").Value = "0.08024"
$ws.Range("E49").Value = "  +1.86%  "
$ws.Range("D50").Value = "131.55"
$ws.Range("E50").Value = "  +0.60%  "
$ws.Range("D51").Value = "1.216"
$ws.Range("E51").Value = "  +2.08%  "

# Remove the temporary text-number-format styling so cells return to their
# original unstyled state (matching the source workbook formatting).
$ws.Range("D2:D51").ClearFormats()
$ws.Range("E2:E51").ClearFormats()
